$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 takes over the old row 19 values (REFREG99 / 99999).
# Copy B19's formatting into B20 first so the numeric style is preserved.
$ws.Range("B19").Copy()
$ws.Range("B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A20").Value = "REFREG99"
$ws.Range("B20").Value = 99999

# Row 19 now becomes REFREG99 / 4000 (was the old row 18 value).
$ws.Range("A19").Value = "REFREG99"
$ws.Range("B19").Value = 4000

# Row 18 becomes the new REFREG93 / 99993 entry.
$ws.Range("A18").Value = "REFREG93"
$ws.Range("B18").Value = 99993
